$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add "ok" status to D2 and D3 (manage game state from outside)
$ws.Range("D2").Value = "ok"
$ws.Range("D3").Value = "ok"

# Update current selection to D4
$ws.Range("D4").Select()
